$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so purely-numeric-looking
# values (e.g. "521.11") are not auto-coerced into Number cells.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "59.354.60"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").Value = "2.649.54"
$ws.Range("E3").Value = "  +0.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").Value = "521.11"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6
$ws.Range("D6").Value = "146.75"
$ws.Range("E6").Value = "  +0.46%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -0.13%  "

# Row 9
$ws.Range("D9").Value = "2.663.37"
$ws.Range("E9").Value = "  +0.47%  "

# Row 10
$ws.Range("E10").Value = "  +9.55%  "

# Row 11
$ws.Range("E11").Value = "  -2.42%  "

# Row 12
$ws.Range("E12").Value = "  -0.46%  "

# Row 13
$ws.Range("E13").Value = "  +1.95%  "

# Row 14
$ws.Range("D14").Value = "3.117.16"
$ws.Range("E14").Value = "  +0.41%  "

# Row 15
$ws.Range("D15").Value = "59.341.67"
$ws.Range("E15").Value = "  +0.12%  "

# Row 16
$ws.Range("E16").Value = "  +0.57%  "

# Row 17
$ws.Range("E17").Value = "  -1.38%  "

# Row 18
$ws.Range("D18").Value = "2.639.91"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19
$ws.Range("D19").Value = "340.29"
$ws.Range("E19").Value = "  -2.75%  "

# Row 20
$ws.Range("E20").Value = "  -1.72%  "

# Row 21
$ws.Range("D21").Value = "10.32"
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$ws.Range("E22").Value = "  +1.32%  "

# Row 23
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").Value = "63.29"
$ws.Range("E24").Value = "  +1.96%  "

# Row 25
$ws.Range("E25").Value = "  +1.47%  "

# Row 26
$ws.Range("E26").Value = "  -1.12%  "

# Row 27
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.32%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0800"
$ws.Range("E28").Value = "  -0.70%  "

# Row 29
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("E30").Value = "  +3.73%  "

# Row 31
$ws.Range("E31").Value = "  +0.08%  "

# Row 32
$ws.Range("E32").Value = "  +0.61%  "

# Row 33
$ws.Range("D33").Value = "18.74"
$ws.Range("E33").Value = "  -1.00%  "

# Row 34
$ws.Range("D34").Value = "149.44"
$ws.Range("E34").Value = "  +0.11%  "

# Row 35
$ws.Range("E35").Value = "  +2.19%  "

# Row 36
$ws.Range("D36").Value = "1.19"
$ws.Range("E36").Value = "  +2.25%  "

# Row 37
$ws.Range("E37").Value = "  -4.73%  "

# Row 38
$ws.Range("D38").Value = "0.881"
$ws.Range("E38").Value = "  +1.34%  "

# Row 39
$ws.Range("D39").Value = "36.92"
$ws.Range("E39").Value = "  +0.85%  "

# Row 40
$ws.Range("D40").Value = "1.48"
$ws.Range("E40").Value = "  +2.26%  "

# Row 41
$ws.Range("E41").Value = "  -2.20%  "

# Row 42
$ws.Range("E42").Value = "  +4.41%  "

# Row 43
$ws.Range("E43").Value = "  +0.24%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.76"
$ws.Range("E44").Value = "  +0.46%  "

# Row 45
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "273.72"
$ws.Range("E45").Value = "  -1.39%  "

# Row 46
$ws.Range("D46").Value = "0.0975"
$ws.Range("E46").Value = "  -1.69%  "

# Row 47
$ws.Range("D47").Value = "0.0535"
$ws.Range("E47").Value = "  +1.31%  "

# Row 48
$ws.Range("D48").Value = "2.050.66"
$ws.Range("E48").Value = "  -2.60%  "

# Row 49
$ws.Range("E49").Value = "  +2.01%  "

# Row 50
$ws.Range("D50").Value = "4.78"
$ws.Range("E50").Value = "  +1.08%  "

# Row 51
$ws.Range("E51").Value = "  -1.18%  "

# Restore the default style on column D so no stray formatting is left
# behind now that the text values are in place.
$priceCol.Style = "Normal"
